$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the stray _GoBack bookmark that sits on the title
# paragraph (it will be re-created later around "DX Grid").
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# Step 2: rewrite the "Click on waterfall ... during the move." bullet
# so it becomes the start of the merged "Double-click on a decoded
# message ..." sentence.
# ---------------------------------------------------------------------
$p50 = $d.Paragraphs.Item(50)
$oldP50 = "Click on waterfall to set QSO Frequency.  Double-clicking will set QSO Frequency, start the decoder, and set Tol to a suitable mode-dependent value.  By default the audio Tx Frequency tracks the QSO Frequency; you can move only the QSO frequency by holding down the CTRL key during the move."
$newP50 = "Double-click on a decoded message to copy the callsign and locator into "
$rngB = $d.Range($p50.Range.Start, $p50.Range.End)
$null = $rngB.Find.Execute($oldP50, $false, $false, $false, $false, $false, $true, 1, $false, $newP50, 2)

# ---------------------------------------------------------------------
# Step 3: strip the now-redundant lead-in text of the next bullet
# ("Double-click on a decoded callsign to copy it into "), leaving the
# bold "DX Call" run untouched.
# ---------------------------------------------------------------------
$p51 = $d.Paragraphs.Item(51)
$rngC = $d.Range($p51.Range.Start, $p51.Range.End)
$oldLead = "Double-click on a decoded callsign to copy it into "
$null = $rngC.Find.Execute($oldLead, $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------
# Step 4: turn " and generate" (immediately after "DX Call") into
# " and DX Grid.  This will also generate", with "DX Grid" bold and
# wrapped in a freshly minted "_GoBack" bookmark.
# ---------------------------------------------------------------------
$p51b = $d.Paragraphs.Item(51)
$rngFindDx = $d.Range($p51b.Range.Start, $p51b.Range.End)
$null = $rngFindDx.Find.Execute("DX Call", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterDxCallStart = $rngFindDx.End

$rngAndGen = $d.Range($afterDxCallStart, $p51b.Range.End)
$null = $rngAndGen.Find.Execute(" and generate", $false, $false, $false, $false, $false, $true, 1, $false, " and DX Grid.  This will also generate", 2)

$p51c = $d.Paragraphs.Item(51)
$rngDxGrid = $d.Range($p51c.Range.Start, $p51c.Range.End)
$null = $rngDxGrid.Find.Execute("DX Grid", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngDxGrid.Font.Bold = -1
$d.Bookmarks.Add("_GoBack", $rngDxGrid)

# ---------------------------------------------------------------------
# Step 5: the two bullets are really one sentence now - merge them by
# deleting the paragraph mark that still separates them.
# ---------------------------------------------------------------------
$p50again = $d.Paragraphs.Item(50)
$markRng = $d.Range($p50again.Range.End - 1, $p50again.Range.End)
$markRng.Delete()

# ---------------------------------------------------------------------
# Step 6: the document lost a page after the rewrite above, so the
# cached " PAGE " field result in the footer needs to go from 6 to 5.
# ---------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$pageChar = $footer.Range.Characters.Item(1)
$pageChar.Text = "5"
